# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp string (06:22 -> 06:52)
# - Refresh Hungria's stats (row 64)
# - Camerun's case count grew enough to overtake Islandia and Azerbaiyan in
#   the ranking, so it moves from row 74 up to row 72; Islandia and
#   Azerbaiyan each shift down one row, keeping their own (unchanged) data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (A1) ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 06:52"

# --- Hungria (row 64) refreshed counts ----------------------------------
$ws.Cells.Item(64, 2).Value = 2727
$ws.Cells.Item(64, 3).Value = 78
$ws.Cells.Item(64, 4).Value = 536
$ws.Cells.Item(64, 5).Value = 1891
$ws.Cells.Item(64, 6).Value = 50
$ws.Cells.Item(64, 7).Value = 9
$ws.Cells.Item(64, 8).Value = 300

# --- Reorder Camerun / Islandia / Azerbaiyan (rows 72-74) ---------------
# Row 72 becomes Camerun with its freshly updated figures.
$ws.Cells.Item(72, 1).Value = "Camerun"
$ws.Cells.Item(72, 2).Value = 1806
$ws.Cells.Item(72, 3).Value = 101
$ws.Cells.Item(72, 4).Value = 915
$ws.Cells.Item(72, 5).Value = 832
$ws.Cells.Item(72, 6).Value = 12
$ws.Cells.Item(72, 7).Value = 1
$ws.Cells.Item(72, 8).Value = 59

# Row 73 becomes Islandia, keeping the figures it had before the reshuffle.
$ws.Cells.Item(73, 1).Value = "Islandia"
$ws.Cells.Item(73, 2).Value = 1795
$ws.Cells.Item(73, 3).Value = 0
$ws.Cells.Item(73, 4).Value = 1636
$ws.Cells.Item(73, 5).Value = 149
$ws.Cells.Item(73, 6).Value = 1
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 10

# Row 74 becomes Azerbaiyan, keeping the figures it had before the reshuffle.
$ws.Cells.Item(74, 1).Value = "Azerbaiyan"
$ws.Cells.Item(74, 2).Value = 1717
$ws.Cells.Item(74, 3).Value = 0
$ws.Cells.Item(74, 4).Value = 1221
$ws.Cells.Item(74, 5).Value = 474
$ws.Cells.Item(74, 6).Value = 15
$ws.Cells.Item(74, 7).Value = 0
$ws.Cells.Item(74, 8).Value = 22
